$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Ensure enough rows exist; insert 2 new rows if sheet currently has 20 data rows (21 total) so it can hold 22 data rows (23 total incl. header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 23) {
    for ($i = $lastRow + 1; $i -le 23; $i++) {
        $ws.Rows.Item($i).Insert()
    }
}

$ws.Cells.Item(2,1).Value = 24330051920092
$ws.Cells.Item(2,2).Value = "APARICIO"
$ws.Cells.Item(2,3).Value = "OFICIAL"
$ws.Cells.Item(2,4).Value = "VICTOR YAEL"
$ws.Cells.Item(2,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(2,6).Value = "2AEV"
$ws.Cells.Item(2,7).Value = 4

$ws.Cells.Item(3,1).Value = 24330051920304
$ws.Cells.Item(3,2).Value = "ARMAS"
$ws.Cells.Item(3,3).Value = "SALINAS"
$ws.Cells.Item(3,4).Value = "JOSE GUSTAVO"
$ws.Cells.Item(3,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(3,6).Value = "2AEV"
$ws.Cells.Item(3,7).Value = 4

$ws.Cells.Item(4,1).Value = 24330051920305
$ws.Cells.Item(4,2).Value = "MORALES"
$ws.Cells.Item(4,3).Value = "CUAHUA"
$ws.Cells.Item(4,4).Value = "ANDRES"
$ws.Cells.Item(4,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(4,6).Value = "2AEV"
$ws.Cells.Item(4,7).Value = 4

$ws.Cells.Item(5,1).Value = 24330051920113
$ws.Cells.Item(5,2).Value = "RAMOS"
$ws.Cells.Item(5,3).Value = "DE LA CRUZ"
$ws.Cells.Item(5,4).Value = "DEREK"
$ws.Cells.Item(5,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(5,6).Value = "2AEV"
$ws.Cells.Item(5,7).Value = 4

$ws.Cells.Item(6,1).Value = 24330051920315
$ws.Cells.Item(6,2).Value = "VENTURA"
$ws.Cells.Item(6,3).Value = "ZEPEDA"
$ws.Cells.Item(6,4).Value = "CARLOS ARGEL"
$ws.Cells.Item(6,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(6,6).Value = "2AEV"
$ws.Cells.Item(6,7).Value = 4

$ws.Cells.Item(7,1).Value = 24330051920330
$ws.Cells.Item(7,2).Value = "VASQUEZ"
$ws.Cells.Item(7,3).Value = "PEREZ"
$ws.Cells.Item(7,4).Value = "DANIELA LILI"
$ws.Cells.Item(7,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(7,6).Value = "2ALCV"
$ws.Cells.Item(7,7).Value = 4

$ws.Cells.Item(8,1).Value = 24330051920246
$ws.Cells.Item(8,2).Value = "ZUNO"
$ws.Cells.Item(8,3).Value = "FLORES"
$ws.Cells.Item(8,4).Value = "ALIN MARIEL"
$ws.Cells.Item(8,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(8,6).Value = "2ALCV"
$ws.Cells.Item(8,7).Value = 4

$ws.Cells.Item(9,1).Value = 24330051920182
$ws.Cells.Item(9,2).Value = "LOPEZ"
$ws.Cells.Item(9,3).Value = "DE LA CRUZ"
$ws.Cells.Item(9,4).Value = "AMISADAY"
$ws.Cells.Item(9,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(9,6).Value = "2ARHV"
$ws.Cells.Item(9,7).Value = 4

$ws.Cells.Item(10,1).Value = 24330051920187
$ws.Cells.Item(10,2).Value = "OSORIO"
$ws.Cells.Item(10,3).Value = "HERNANDEZ"
$ws.Cells.Item(10,4).Value = "AYLIN ABIGAIL"
$ws.Cells.Item(10,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(10,6).Value = "2ARHV"
$ws.Cells.Item(10,7).Value = 4

$ws.Cells.Item(11,1).Value = 24330051920093
$ws.Cells.Item(11,2).Value = "ARIAS"
$ws.Cells.Item(11,3).Value = "SARMIENTO"
$ws.Cells.Item(11,4).Value = "URIEL ARTURO"
$ws.Cells.Item(11,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(11,6).Value = "2AEV"
$ws.Cells.Item(11,7).Value = 3

$ws.Cells.Item(12,1).Value = 24330051920090
$ws.Cells.Item(12,2).Value = "ANTONIO"
$ws.Cells.Item(12,3).Value = "LOPEZ"
$ws.Cells.Item(12,4).Value = "SERGIO GISELL"
$ws.Cells.Item(12,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(12,6).Value = "2AEV"
$ws.Cells.Item(12,7).Value = 3

$ws.Cells.Item(13,1).Value = 24330051920098
$ws.Cells.Item(13,2).Value = "CHICO"
$ws.Cells.Item(13,3).Value = "BALDERAS"
$ws.Cells.Item(13,4).Value = "YARETH"
$ws.Cells.Item(13,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(13,6).Value = "2AEV"
$ws.Cells.Item(13,7).Value = 3

$ws.Cells.Item(14,1).Value = 23330051920036
$ws.Cells.Item(14,2).Value = "HERNANDEZ"
$ws.Cells.Item(14,3).Value = "DOLORES"
$ws.Cells.Item(14,4).Value = "GONZALO"
$ws.Cells.Item(14,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(14,6).Value = "2AEV"
$ws.Cells.Item(14,7).Value = 3

$ws.Cells.Item(15,1).Value = 24330051920144
$ws.Cells.Item(15,2).Value = "MUÑOZ"
$ws.Cells.Item(15,3).Value = "CORONA"
$ws.Cells.Item(15,4).Value = "JOSE ABEL"
$ws.Cells.Item(15,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(15,6).Value = "2AEV"
$ws.Cells.Item(15,7).Value = 3

$ws.Cells.Item(16,1).Value = 24330051920306
$ws.Cells.Item(16,2).Value = "ROJAS"
$ws.Cells.Item(16,3).Value = "GUTIERREZ"
$ws.Cells.Item(16,4).Value = "LUIS ROBERTO"
$ws.Cells.Item(16,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(16,6).Value = "2AEV"
$ws.Cells.Item(16,7).Value = 3

$ws.Cells.Item(17,1).Value = 24330051920398
$ws.Cells.Item(17,2).Value = "CORTES"
$ws.Cells.Item(17,3).Value = "MENDEZ"
$ws.Cells.Item(17,4).Value = "ELIZABETH MADAI"
$ws.Cells.Item(17,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(17,6).Value = "2ALCV"
$ws.Cells.Item(17,7).Value = 3

$ws.Cells.Item(18,1).Value = 24330051920389
$ws.Cells.Item(18,2).Value = "RUIZ"
$ws.Cells.Item(18,3).Value = "MORALES"
$ws.Cells.Item(18,4).Value = "MAYRIN GUADALUPE"
$ws.Cells.Item(18,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(18,6).Value = "2ALCV"
$ws.Cells.Item(18,7).Value = 3

$ws.Cells.Item(19,1).Value = 23330051920298
$ws.Cells.Item(19,2).Value = "MAZA"
$ws.Cells.Item(19,3).Value = "ENCARNACION"
$ws.Cells.Item(19,4).Value = "KEVIN JESUS"
$ws.Cells.Item(19,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(19,6).Value = "2ARHV"
$ws.Cells.Item(19,7).Value = 3

$ws.Cells.Item(20,1).Value = 23330051920224
$ws.Cells.Item(20,2).Value = "DORANTES"
$ws.Cells.Item(20,3).Value = "PORRAS"
$ws.Cells.Item(20,4).Value = "ROBERTO"
$ws.Cells.Item(20,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(20,6).Value = "2AEV"
$ws.Cells.Item(20,7).Value = 2

$ws.Cells.Item(21,1).Value = 24330051920238
$ws.Cells.Item(21,2).Value = "TORRES"
$ws.Cells.Item(21,3).Value = "PEREZ"
$ws.Cells.Item(21,4).Value = "ERIKA VALERIA"
$ws.Cells.Item(21,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(21,6).Value = "2ALCV"
$ws.Cells.Item(21,7).Value = 2

$ws.Cells.Item(22,1).Value = 24330051920206
$ws.Cells.Item(22,2).Value = "PORTUGAL"
$ws.Cells.Item(22,3).Value = "VEGA"
$ws.Cells.Item(22,4).Value = "SANTIAGO"
$ws.Cells.Item(22,5).Value = "Pensamiento matemático II"
$ws.Cells.Item(22,6).Value = "2ARHV"
$ws.Cells.Item(22,7).Value = 2
